$d = $word.ActiveDocument

# 1) "Top5-list" -> "Top10-list"
$d.Content.Find.Execute("Top5-list", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Top10-list", 2) | Out-Null

# 2) "all time best" -> "all-time best"
$d.Content.Find.Execute("all time best", $true, $false, $false, $false, $false,
                         $true, 1, $false, "all-time best", 2) | Out-Null
